$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'31.217.70"
$ws.Range("E2").Value = "  +2.38%  "
$ws.Range("D3").Value = "'1.972.82"
$ws.Range("E3").Value = "  +3.31%  "
$ws.Range("D4").Value = "'1.004"
$ws.Range("E4").Value = "  +0.33%  "
$ws.Range("D5").Value = "'248.83"
$ws.Range("E5").Value = "  +2.02%  "
$ws.Range("E6").Value = "  +0.31%  "
$ws.Range("D7").Value = "'0.4901"
$ws.Range("E7").Value = "  +1.32%  "
$ws.Range("D8").Value = "'44.90"
$ws.Range("E8").Value = "  +1.05%  "
$ws.Range("D9").Value = "'0.2959"
$ws.Range("E9").Value = "  +2.75%  "
$ws.Range("D10").Value = "'0.06846"
$ws.Range("E10").Value = "  +0.61%  "
$ws.Range("D11").Value = "'19.23"
$ws.Range("E11").Value = "  -0.66%  "
$ws.Range("D12").Value = "'107.58"
$ws.Range("E12").Value = "  -3.28%  "
$ws.Range("D13").Value = "'1.967.61"
$ws.Range("E13").Value = "  +2.80%  "
$ws.Range("D14").Value = "'0.07795"
$ws.Range("E14").Value = "  +3.13%  "
$ws.Range("D15").Value = "'5.455"
$ws.Range("E15").Value = "  +1.33%  "
$ws.Range("D16").Value = "'0.7102"
$ws.Range("E16").Value = "  +6.39%  "
$ws.Range("D17").Value = "'285.96"
$ws.Range("E17").Value = "  -2.40%  "
$ws.Range("D18").Value = "'31.206.82"
$ws.Range("E18").Value = "  +2.35%  "
$ws.Range("D19").Value = "'13.32"
$ws.Range("E19").Value = "  +2.57%  "
$ws.Range("D20").Value = "'0.000007758"
$ws.Range("E20").Value = "  +2.71%  "
$ws.Range("D21").Value = "'2.222.44"
$ws.Range("E21").Value = "  +2.91%  "
$ws.Range("D22").Value = "'5.634"
$ws.Range("E22").Value = "  +2.65%  "
$ws.Range("E23").Value = "  +0.42%  "
$ws.Range("D24").Value = "'1.005"
$ws.Range("E24").Value = "  +0.43%  "
$ws.Range("D25").Value = "'6.661"
$ws.Range("E25").Value = "  +4.45%  "
$ws.Range("D26").Value = "'10.05"
$ws.Range("E26").Value = "  +6.48%  "
$ws.Range("D27").Value = "'170.24"
$ws.Range("E27").Value = "  +3.01%  "
$ws.Range("D28").Value = "'20.17"
$ws.Range("E28").Value = "  -0.25%  "
$ws.Range("D29").Value = "'2.198"
$ws.Range("E29").Value = "  +6.07%  "
$ws.Range("E30").Value = "  +0.31%  "
$ws.Range("D31").Value = "'1.448"
$ws.Range("E31").Value = "  +1.34%  "
$ws.Range("D32").Value = "'4.858"
$ws.Range("E32").Value = "  +20.00%  "
$ws.Range("E33").Value = "  +9.66%  "
$ws.Range("D34").Value = "'0.05090"
$ws.Range("E34").Value = "  +2.40%  "
$ws.Range("D35").Value = "'0.7731"
$ws.Range("E35").Value = "  +5.34%  "
$ws.Range("D36").Value = "'1.174"
$ws.Range("E36").Value = "  +3.86%  "
$ws.Range("D37").Value = "'2.745"
$ws.Range("E37").Value = "  +0.99%  "
$ws.Range("D38").Value = "'0.02052"
$ws.Range("E38").Value = "  +0.97%  "
$ws.Range("D39").Value = "'2.738"
$ws.Range("E39").Value = "  +2.05%  "
$ws.Range("D40").Value = "'6.463"
$ws.Range("E40").Value = "  +11.97%  "
$ws.Range("D41").Value = "'2.136"
$ws.Range("E41").Value = "  +6.52%  "
$ws.Range("D42").Value = "'74.41"
$ws.Range("E42").Value = "  +7.98%  "
$ws.Range("D43").Value = "'0.8884"
$ws.Range("E43").Value = "  +2.96%  "
$ws.Range("D44").Value = "'110.21"
$ws.Range("E44").Value = "  +0.90%  "
$ws.Range("D45").Value = "'0.4482"
$ws.Range("E45").Value = "  +1.51%  "
$ws.Range("E46").Value = "  +0.42%  "
$ws.Range("D47").Value = "'7.533"
$ws.Range("E47").Value = "  +4.98%  "
$ws.Range("D48").Value = "'989.68"
$ws.Range("E48").Value = "  +17.40%  "
$ws.Range("D49").Value = "'9.465"
$ws.Range("E49").Value = "  +2.64%  "
$ws.Range("D50").Value = "'0.1276"
$ws.Range("E50").Value = "  +4.15%  "
$ws.Range("D51").Value = "'36.09"
$ws.Range("E51").Value = "  +4.15%  "
